# إضافة حدث جديد في Card1
# - fills the previously-empty "nan" placeholder cells B22:K22
# - appends a new event row (row 23) with date/event/correction/servicedBy
# - worksheet dimension grows from A1:O22 to A1:O23 automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# --- Row 22: the "Min_Tones..Revolving flats(o)" columns (B:K) were left
#     as blank placeholders; they should now hold the literal text "nan"
#     (matching the rest of the sheet's convention for missing numeric data).
$emptyCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $emptyCols) {
    $ws.Range($col + "22").Value = "nan"
}

# --- Row 23: brand-new event row.
# Column A holds the text "1" (like every other row on this sheet). Excel
# would normally infer a number from a bare "1", so use a leading
# apostrophe to force text, then clear the resulting quote-prefix style so
# no extra formatting is introduced.
$ws.Range("A23").Value = "'1"
$ws.Range("A23").Style = "Normal"

# Columns B:K stay blank for the new row, same as the other data rows.
# (briefly give the cell a format so it is materialised in the sheet even
# once its value is cleared back to empty, then restore the default style)
foreach ($col in $emptyCols) {
    $cell = $ws.Range($col + "23")
    $cell.Value = "x"
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.Style = "Normal"
}

$ws.Range("L23").Value = "21/2/2026"
$ws.Range("M23").Value = "قطع سير700"
$ws.Range("N23").Value = "تم تغير سير 700(مشلان)"
$ws.Range("O23").Value = "رضا"
